$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E11").Value = 113880
